# Applies the cryptos.xlsx price/volume refresh described in the commit:
#   "Updated cryptos list on Sat Sep 28 15:18:33 UTC 2024 with GitHub Actions"
#
# Rows 2-51 (Bitcoin..EnergySwap): Price (D) and Volume(1h) (E) columns are
# refreshed with newly scraped values. Rows 38/39 additionally swap their
# rank order (Monero <-> FirstDigitalUSD), so Coin (B) and Link (C) change too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell reference + its new text.
# Column D holds numeric-looking strings (prices such as '602.91' or
# thousand-dotted '65.712.78'); Excel's COM layer auto-converts plain
# numeric-looking assignments to a Double, which would both lose the
# original text formatting (e.g. trailing zeros) and change the stored
# cell type from string to number. Forcing NumberFormat to "@" (Text)
# before the assignment keeps the value a literal string, matching the
# source workbook; resetting the style back to "Normal" afterwards avoids
# leaving a stray text-format style on the cell.
$updates = @(
    @{Ref="D2"; Value="65.712.78"}
    @{Ref="E2"; Value="  -0.97%  "}
    @{Ref="D3"; Value="2.671.73"}
    @{Ref="E3"; Value="  -0.66%  "}
    @{Ref="E4"; Value="  +0.04%  "}
    @{Ref="D5"; Value="602.91"}
    @{Ref="E5"; Value="  -1.39%  "}
    @{Ref="D6"; Value="157.63"}
    @{Ref="E6"; Value="  -2.00%  "}
    @{Ref="E7"; Value="  +0.05%  "}
    @{Ref="D8"; Value="0.620"}
    @{Ref="E8"; Value="  +4.57%  "}
    @{Ref="D9"; Value="0.131"}
    @{Ref="E9"; Value="  +3.47%  "}
    @{Ref="E10"; Value="  -1.02%  "}
    @{Ref="E11"; Value="  -3.41%  "}
    @{Ref="E12"; Value="  -0.31%  "}
    @{Ref="D13"; Value="29.58"}
    @{Ref="E13"; Value="  -2.70%  "}
    @{Ref="E14"; Value="  -5.67%  "}
    @{Ref="D15"; Value="3.154.09"}
    @{Ref="E15"; Value="  -0.55%  "}
    @{Ref="D16"; Value="65.547.09"}
    @{Ref="E16"; Value="  -0.88%  "}
    @{Ref="D17"; Value="2.660.98"}
    @{Ref="E17"; Value="  -0.91%  "}
    @{Ref="D18"; Value="12.86"}
    @{Ref="E18"; Value="  +0.86%  "}
    @{Ref="E19"; Value="  -2.03%  "}
    @{Ref="D20"; Value="7.68"}
    @{Ref="E20"; Value="  +2.30%  "}
    @{Ref="D21"; Value="352.09"}
    @{Ref="E21"; Value="  -3.26%  "}
    @{Ref="D22"; Value="1.00"}
    @{Ref="E22"; Value="  -0.13%  "}
    @{Ref="D23"; Value="69.92"}
    @{Ref="E23"; Value="  -0.42%  "}
    @{Ref="E24"; Value="  +2.79%  "}
    @{Ref="D25"; Value="9.78"}
    @{Ref="E25"; Value="  +0.21%  "}
    @{Ref="E26"; Value="  -1.68%  "}
    @{Ref="D27"; Value="0.168"}
    @{Ref="E27"; Value="  -4.08%  "}
    @{Ref="E28"; Value="  -5.39%  "}
    @{Ref="E29"; Value="  -1.20%  "}
    @{Ref="E30"; Value="  -0.12%  "}
    @{Ref="D31"; Value="2.14"}
    @{Ref="E31"; Value="  -3.19%  "}
    @{Ref="D32"; Value="530.62"}
    @{Ref="E32"; Value="  -2.60%  "}
    @{Ref="D33"; Value="1.76"}
    @{Ref="E33"; Value="  -2.61%  "}
    @{Ref="E34"; Value="  -1.32%  "}
    @{Ref="D35"; Value="5.53"}
    @{Ref="E35"; Value="  +0.62%  "}
    @{Ref="D36"; Value="0.425"}
    @{Ref="E36"; Value="  -2.85%  "}
    @{Ref="D37"; Value="20.48"}
    @{Ref="E37"; Value="  -1.74%  "}
    @{Ref="B38"; Value="FirstDigitalUSD"}
    @{Ref="C38"; Value="https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"}
    @{Ref="D38"; Value="0.999"}
    @{Ref="E38"; Value="  +0.01%  "}
    @{Ref="B39"; Value="Monero"}
    @{Ref="C39"; Value="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"}
    @{Ref="D39"; Value="159.21"}
    @{Ref="E39"; Value="  -2.34%  "}
    @{Ref="E40"; Value="  -3.21%  "}
    @{Ref="D41"; Value="0.999"}
    @{Ref="E41"; Value="  -0.02%  "}
    @{Ref="D42"; Value="42.57"}
    @{Ref="E42"; Value="  +0.01%  "}
    @{Ref="D43"; Value="165.66"}
    @{Ref="E43"; Value="  -2.90%  "}
    @{Ref="E44"; Value="  -3.32%  "}
    @{Ref="D45"; Value="2.32"}
    @{Ref="E45"; Value="  -0.95%  "}
    @{Ref="E46"; Value="  -1.40%  "}
    @{Ref="D47"; Value="23.18"}
    @{Ref="E47"; Value="  -0.57%  "}
    @{Ref="E48"; Value="  -3.09%  "}
    @{Ref="E49"; Value="  -3.03%  "}
    @{Ref="D50"; Value="0.101"}
    @{Ref="E50"; Value="  +2.30%  "}
    @{Ref="D51"; Value="20.20"}
    @{Ref="E51"; Value="  -0.55%  "}
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Ref)
    if ($u.Ref.StartsWith("D")) {
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}

Write-Host "Applied $($updates.Count) cell updates."
